$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 144, shifting the existing rows 144:211 down to 145:212.
$ws.Rows(144).Insert()

# Populate the newly inserted row 144 with the new weekly data point.
$ws.Cells.Item(144, 1).Value = 9
$ws.Cells.Item(144, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(144, 3).Value = 'Metropolitana'
$ws.Cells.Item(144, 4).Value = 44572
$ws.Cells.Item(144, 5).Value = 13
$ws.Cells.Item(144, 6).Value = 100112043
$ws.Cells.Item(144, 7).Value = 'Pepino ensalada'
$ws.Cells.Item(144, 8).Value = 'Sin especificar'
$ws.Cells.Item(144, 9).Value = 'Primera'
$ws.Cells.Item(144, 10).Value = 160
$ws.Cells.Item(144, 11).Value = 6000
$ws.Cells.Item(144, 12).Value = 7000
$ws.Cells.Item(144, 13).Value = 6500
$ws.Cells.Item(144, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(144, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(144, 16).Value = 108
$ws.Cells.Item(144, 17).Value = 60
$ws.Cells.Item(144, 18).Value = 'Hortaliza'
